{"js": "const pairs = [\n  [\"2024-10-25 Friday\", \"2024-10-26 Saturday\"],\n  [\"59-2=57\", \"66-9=57\"],\n  [\"72+18=90\", \"31-3=28\"],\n  [\"57-49=8\", \"51+7=58\"],\n  [\"73-25=48\", \"15+33=48\"],\n  [\"47-30=17\", \"44+0=44\"],\n  [\"64-21=43\", \"55+7=62\"],\n  [\"62-53=9\", \"15+70=85\"],\n  [\"73+9=82\", \"81-64=17\"],\n  [\"98-15=83\", \"45+23=68\"],\n  [\"11+25=36\", \"71+21=92\"],\n  [\"73+18=91\", \"37+18=55\"],\n  [\"39-38=1\", \"91-46=45\"],\n  [\"59-3=56\", \"34+3=37\"],\n  [\"74-56=18\", \"0+24=24\"],\n  [\"9+24=33\", \"18+8=26\"],\n  [\"71-10=61\", \"41+7=48\"],\n  [\"21+11=32\", \"50+31=81\"],\n  [\"31+16=47\", \"81-54=27\"],\n  [\"51-31=20\", \"54-19=35\"],\n  [\"8+78=86\", \"78-63=15\"],\n  [\"7+61=68\", \"86-28=58\"],\n  [\"52-4=48\", \"78-54=24\"],\n  [\"33-18=15\", \"12+55=67\"],\n  [\"39+20=59\", \"3+54=57\"],\n  [\"74-69=5\", \"20-7=13\"],\n  [\"54+13=67\", \"57+1=58\"],\n  [\"21+68=89\", \"79-37=42\"],\n  [\"35+44=79\", \"49-2=47\"],\n  [\"60-13=47\", \"98-64=34\"],\n  [\"50-42=8\", \"22+37=59\"],\n  [\"9+29=38\", \"56+42=98\"],\n  [\"13+8=21\", \"53-27=26\"],\n  [\"37-22=15\", \"53+4=57\"],\n  [\"54+5=59\", \"93-14=79\"],\n  [\"92-36=56\", \"83-57=26\"],\n  [\"78+3=81\", \"41+26=67\"],\n  [\"77-3=74\", \"31-10=21\"],\n  [\"17+8=25\", \"85-81=4\"],\n  [\"40-23=17\", \"19-2=17\"],\n  [\"89-63=26\", \"1+95=96\"],\n  [\"9+26=35\", \"82-47=35\"],\n  [\"68+18=86\", \"1+47=48\"],\n  [\"65+11=76\", \"97-48=49\"],\n  [\"31+41=72\", \"90+4=94\"],\n  [\"27+69=96\", \"50+24=74\"],\n  [\"26+5=31\", \"43-38=5\"],\n  [\"18-11=7\", \"37+22=59\"],\n  [\"6+71=77\", \"25+62=87\"],\n  [\"41-22=19\", \"39-0=39\"],\n  [\"86+8=94\", \"62-55=7\"],\n  [\"68-0=68\", \"0+62=62\"],\n  [\"50-3=47\", \"1+91=92\"],\n  [\"6+1=7\", \"69-10=59\"],\n  [\"26+16=42\", \"45+49=94\"],\n  [\"32+42=74\", \"27+32=59\"],\n  [\"11-3=8\", \"9+87=96\"],\n  [\"30+33=63\", \"87-5=82\"],\n  [\"8+64=72\", \"50-30=20\"],\n  [\"11+5=16\", \"54-37=17\"],\n  [\"14+27=41\", \"5-3=2\"],\n  [\"80+7=87\", \"70-5=65\"],\n  [\"39-37=2\", \"20-3=17\"],\n  [\"14+61=75\", \"16-3=13\"],\n  [\"30-18=12\", \"70-45=25\"],\n  [\"45-26=19\", \"12+79=91\"],\n  [\"39-28=11\", \"8+29=37\"],\n  [\"80-3=77\", \"93-16=77\"],\n  [\"41-31=10\", \"85-7=78\"],\n  [\"3+20=23\", \"97-85=12\"],\n  [\"32+39=71\", \"71-70=1\"],\n  [\"28+4=32\", \"26+11=37\"],\n  [\"47+24=71\", \"94-57=37\"],\n  [\"75-61=14\", \"57-1=56\"],\n  [\"40+49=89\", \"23-19=4\"],\n  [\"91-26=65\", \"49-35=14\"],\n  [\"2+90=92\", \"54+26=80\"],\n  [\"98-80=18\", \"9+59=68\"],\n  [\"11+44=55\", \"41-21=20\"],\n  [\"92-24=68\", \"28+13=41\"],\n  [\"65-40=25\", \"81-21=60\"],\n  [\"1+94=95\", \"15+22=37\"],\n  [\"87+7=94\", \"49+44=93\"],\n  [\"38-30=8\", \"36+27=63\"],\n  [\"97-86=11\", \"27+43=70\"],\n  [\"92-37=55\", \"51-14=37\"],\n  [\"21+22=43\", \"27+36=63\"],\n  [\"28+15=43\", \"32+33=65\"],\n  [\"99-51=48\", \"40+38=78\"],\n  [\"40-37=3\", \"68-66=2\"],\n  [\"18-15=3\", \"53-3=50\"],\n  [\"54+8=62\", \"4+24=28\"],\n  [\"28-20=8\", \"93+3=96\"],\n  [\"68-59=9\", \"49+5=54\"],\n  [\"36+6=42\", \"38+39=77\"],\n  [\"65+3=68\", \"42-3=39\"],\n  [\"68+10=78\", \"32-28=4\"],\n  [\"17+56=73\", \"89-69=20\"],\n  [\"14+21=35\", \"18+31=49\"],\n  [\"29+26=55\", \"89-10=79\"],\n  [\"56+11=67\", \"87-49=38\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load('items');\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@('2024-10-25 Friday', '2024-10-26 Saturday')\n  ,@('59-2=57', '66-9=57')\n  ,@('72+18=90', '31-3=28')\n  ,@('57-49=8', '51+7=58')\n  ,@('73-25=48', '15+33=48')\n  ,@('47-30=17', '44+0=44')\n  ,@('64-21=43', '55+7=62')\n  ,@('62-53=9', '15+70=85')\n  ,@('73+9=82', '81-64=17')\n  ,@('98-15=83', '45+23=68')\n  ,@('11+25=36', '71+21=92')\n  ,@('73+18=91', '37+18=55')\n  ,@('39-38=1', '91-46=45')\n  ,@('59-3=56', '34+3=37')\n  ,@('74-56=18', '0+24=24')\n  ,@('9+24=33', '18+8=26')\n  ,@('71-10=61', '41+7=48')\n  ,@('21+11=32', '50+31=81')\n  ,@('31+16=47', '81-54=27')\n  ,@('51-31=20', '54-19=35')\n  ,@('8+78=86', '78-63=15')\n  ,@('7+61=68', '86-28=58')\n  ,@('52-4=48', '78-54=24')\n  ,@('33-18=15', '12+55=67')\n  ,@('39+20=59', '3+54=57')\n  ,@('74-69=5', '20-7=13')\n  ,@('54+13=67', '57+1=58')\n  ,@('21+68=89', '79-37=42')\n  ,@('35+44=79', '49-2=47')\n  ,@('60-13=47', '98-64=34')\n  ,@('50-42=8', '22+37=59')\n  ,@('9+29=38', '56+42=98')\n  ,@('13+8=21', '53-27=26')\n  ,@('37-22=15', '53+4=57')\n  ,@('54+5=59', '93-14=79')\n  ,@('92-36=56', '83-57=26')\n  ,@('78+3=81', '41+26=67')\n  ,@('77-3=74', '31-10=21')\n  ,@('17+8=25', '85-81=4')\n  ,@('40-23=17', '19-2=17')\n  ,@('89-63=26', '1+95=96')\n  ,@('9+26=35', '82-47=35')\n  ,@('68+18=86', '1+47=48')\n  ,@('65+11=76', '97-48=49')\n  ,@('31+41=72', '90+4=94')\n  ,@('27+69=96', '50+24=74')\n  ,@('26+5=31', '43-38=5')\n  ,@('18-11=7', '37+22=59')\n  ,@('6+71=77', '25+62=87')\n  ,@('41-22=19', '39-0=39')\n  ,@('86+8=94', '62-55=7')\n  ,@('68-0=68', '0+62=62')\n  ,@('50-3=47', '1+91=92')\n  ,@('6+1=7', '69-10=59')\n  ,@('26+16=42', '45+49=94')\n  ,@('32+42=74', '27+32=59')\n  ,@('11-3=8', '9+87=96')\n  ,@('30+33=63', '87-5=82')\n  ,@('8+64=72', '50-30=20')\n  ,@('11+5=16', '54-37=17')\n  ,@('14+27=41', '5-3=2')\n  ,@('80+7=87', '70-5=65')\n  ,@('39-37=2', '20-3=17')\n  ,@('14+61=75', '16-3=13')\n  ,@('30-18=12', '70-45=25')\n  ,@('45-26=19', '12+79=91')\n  ,@('39-28=11', '8+29=37')\n  ,@('80-3=77', '93-16=77')\n  ,@('41-31=10', '85-7=78')\n  ,@('3+20=23', '97-85=12')\n  ,@('32+39=71', '71-70=1')\n  ,@('28+4=32', '26+11=37')\n  ,@('47+24=71', '94-57=37')\n  ,@('75-61=14', '57-1=56')\n  ,@('40+49=89', '23-19=4')\n  ,@('91-26=65', '49-35=14')\n  ,@('2+90=92', '54+26=80')\n  ,@('98-80=18', '9+59=68')\n  ,@('11+44=55', '41-21=20')\n  ,@('92-24=68', '28+13=41')\n  ,@('65-40=25', '81-21=60')\n  ,@('1+94=95', '15+22=37')\n  ,@('87+7=94', '49+44=93')\n  ,@('38-30=8', '36+27=63')\n  ,@('97-86=11', '27+43=70')\n  ,@('92-37=55', '51-14=37')\n  ,@('21+22=43', '27+36=63')\n  ,@('28+15=43', '32+33=65')\n  ,@('99-51=48', '40+38=78')\n  ,@('40-37=3', '68-66=2')\n  ,@('18-15=3', '53-3=50')\n  ,@('54+8=62', '4+24=28')\n  ,@('28-20=8', '93+3=96')\n  ,@('68-59=9', '49+5=54')\n  ,@('36+6=42', '38+39=77')\n  ,@('65+3=68', '42-3=39')\n  ,@('68+10=78', '32-28=4')\n  ,@('17+56=73', '89-69=20')\n  ,@('14+21=35', '18+31=49')\n  ,@('29+26=55', '89-10=79')\n  ,@('56+11=67', '87-49=38')\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $range.Find.Text = $old\n  $range.Find.Replacement.Text = $new\n  $range.Find.Forward = $true\n  $range.Find.Wrap = 1\n  $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
